# Applies the "Reorgs, devices designations, and stock part selection" edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SMD-PB")
$ws2 = $wb.Worksheets.Item("TH-PB")
$ws1.Name = "Buttons-SMD"
$ws2.Name = "Buttons-TH"

# ---------------------------------------------------------------------------
# 2. Sheet "Buttons-SMD": remove the two blank separator rows (7 & 8), which
#    shifts the old rows 9/10 up to become rows 7/8.
# ---------------------------------------------------------------------------
$ws1.Range("A7:A8").EntireRow.Delete()

# Tag rows 2-5 (FSM*JSMA parts) and the two shifted rows (7,8) with their
# Stock / Device designation.
foreach ($r in 2,3,4,5) {
    $ws1.Cells.Item($r, 29).Value = "Stock"        # column AC
    $ws1.Cells.Item($r, 32).Value = "SMD-VERTICAL"  # column AF
}
$ws1.Cells.Item(7, 29).Value = "Stock"
$ws1.Cells.Item(7, 32).Value = "SMD-RIGHT-ANGLE2"
$ws1.Cells.Item(8, 29).Value = "Stock"
$ws1.Cells.Item(8, 32).Value = "SMD-RIGHT-ANGLE"

# New row 9: C&K Components PTS525 SMD part.
$ws1.Cells.Item(9, 1).Value = "http://media.digikey.com/pdf/Data%20Sheets/C&K/PTS525_Series_RevJul_2012.pdf"
$ws1.Cells.Item(9, 2).Value = "http://media.digikey.com/photos/CK%20Comp%20Photos/PTS525SM10SMTR_sml.jpg"
$ws1.Cells.Item(9, 3).Value = "CKN9104CT-NDP"
$ws1.Cells.Item(9, 4).Value = "TS525SM15SMTR2"
$ws1.Cells.Item(9, 5).Value = "LFS C&K Components"
$ws1.Cells.Item(9, 6).Value = "SWITCH TACTILE SPST-NO 0.05A 12V"
$ws1.Cells.Item(9, 7).Value = 2726
$ws1.Cells.Item(9, 8).Value = 0
$ws1.Cells.Item(9, 9).Value = 0.65
$ws1.Cells.Item(9, 10).Value = 0
$ws1.Cells.Item(9, 11).Value = 1
$ws1.Cells.Item(9, 12).Value = "Cut Tape (CT)"
$ws1.Cells.Item(9, 13).Value = "PTS525"
$ws1.Cells.Item(9, 14).Value = "SPST-NO"
$ws1.Cells.Item(9, 15).Value = "Off-Mom"
$ws1.Cells.Item(9, 16).Value = "0.05A @ 12VDC"
$ws1.Cells.Item(9, 17).Value = "Standard"
$ws1.Cells.Item(9, 18).Value = "Surface Mount"
$ws1.Cells.Item(9, 19).Value = "1.50mm"
$ws1.Cells.Item(9, 20).Value = "Top Actuated"
$ws1.Cells.Item(9, 21).Value = "Gull Wing"
$ws1.Cells.Item(9, 22).Value = "5.25mm x 5.25mm"
$ws1.Cells.Item(9, 23).Value = "Non-Illuminated"
$ws1.Cells.Item(9, 24).Value = "-"
$ws1.Cells.Item(9, 25).Value = "-"
$ws1.Cells.Item(9, 26).Value = "160gf"
$ws1.Cells.Item(9, 27).Value = "-"
$ws1.Cells.Item(9, 28).Value = "-"

# ---------------------------------------------------------------------------
# 3. Sheet "Buttons-TH": remove the two blank separator rows (10 & 11), which
#    shifts the old rows 12-16 up to become rows 10-14.
# ---------------------------------------------------------------------------
$ws2.Range("A10:A11").EntireRow.Delete()

foreach ($r in 2,3,4,5,6,7,8,9) {
    $ws2.Cells.Item($r, 29).Value = "Stock"        # column AC
    $ws2.Cells.Item($r, 32).Value = "TH-VERTICAL"   # column AF
}
# Row 9's Package cell was previously an empty styled placeholder; populate it.
$ws2.Cells.Item(9, 31).Value = "TE_x-1825910-y"    # column AE

foreach ($r in 10,11,12,13,14) {
    $ws2.Cells.Item($r, 29).Value = "Stock"        # column AC
    $ws2.Cells.Item($r, 32).Value = "TH-VERTICAL"   # column AF
}

# The Digi-Key part number in (now) row 12 was re-typed in black text.
$ws2.Range("C12").Font.Color = 0
